$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row for the item "حنه جلوري سوده 1 كيس" (row 42).
# This shifts all rows below it up by one, updates merged cells, and
# removes the now-unused shared string for that item name.
$ws.Rows.Item(42).Delete()

# Column A holds a simple sequential row counter ("م" = number) that is
# stored as a literal value, so it must be renumbered after the shift.
for ($r = 42; $r -le 50; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 3
}

# The total in column K (now row 51) is a literal sum, not a formula, so
# it must be reduced by the deleted row's price (45) by hand.
$ws.Cells.Item(51, 11).Value = 2327.71

# Row heights in this sheet are fixed per-row (not tied to the shifted
# content), so restore each row's own original height after the shift.
$rowHeights = @{
    42 = 25.5
    43 = 25.5
    44 = 24.75
    45 = 25.5
    46 = 24.75
    47 = 25.5
    48 = 25.5
    49 = 24.75
    50 = 25.5
    51 = 25.5
    52 = 16.5
}
foreach ($r in $rowHeights.Keys) {
    $ws.Rows.Item($r).RowHeight = $rowHeights[$r]
}
